# Update the "Förändrad" (Changed) date column C for rows 2-10 from
# 2023-11-13 (serial 45243) to 2023-11-14 (serial 45244).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value = 45244
    }
}
